# Applies the "2nd testcase and elective basket scheduling with common slots" edit:
#  1. Section_A (sheet1) and Section_B (sheet2): replace the cells that used to hold
#     the courses that are being pulled out of the fixed grid with "Free", and mark the
#     two common elective slots with "<CODE> (Elective)".
#  2. Add a new "Course_Summary" sheet describing the elective basket.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Section_A
# ---------------------------------------------------------------------------
$secA = $wb.Worksheets.Item("Section_A")

$secA.Range("E2").Value = "Free"
$secA.Range("F2").Value = "Free"

$secA.Range("B3").Value = "Free"
$secA.Range("C3").Value = "Free"
$secA.Range("D3").Value = "Free"
$secA.Range("E3").Value = "Free"

$secA.Range("D5").Value = "Free"
$secA.Range("F5").Value = "EC456 (Elective)"

$secA.Range("B6").Value = "Free"
$secA.Range("E6").Value = "DS456 (Elective)"

$secA.Range("B7").Value = "Free"
$secA.Range("C7").Value = "Free"
$secA.Range("D7").Value = "Free"
$secA.Range("E7").Value = "Free"

# ---------------------------------------------------------------------------
# 2. Section_B
# ---------------------------------------------------------------------------
$secB = $wb.Worksheets.Item("Section_B")

$secB.Range("B2").Value = "Free"
$secB.Range("C2").Value = "Free"
$secB.Range("F2").Value = "Free"

$secB.Range("B3").Value = "Free"
$secB.Range("E3").Value = "Free"

$secB.Range("D5").Value = "Free"
$secB.Range("F5").Value = "EC456 (Elective)"

$secB.Range("C6").Value = "Free"
$secB.Range("E6").Value = "DS456 (Elective)"

$secB.Range("B7").Value = "Free"
$secB.Range("C7").Value = "Free"
$secB.Range("D7").Value = "Free"
$secB.Range("F7").Value = "Free"

# ---------------------------------------------------------------------------
# 3. New sheet: Course_Summary (added after the last sheet, i.e. Section_B)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$summary = $wb.Worksheets.Add($null, $lastSheet)
$summary.Name = "Course_Summary"

$headers = @("Course Code", "Course Name", "Course Type", "LTPSC", "Credits", "Instructor")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $summary.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Reuse the bold/bordered header style already used by the timetable sheets
# instead of building a fresh one.
$secA.Range("B1").Copy()
$summary.Range("A1:F1").PasteSpecial(-4122)

$rows = @(
    @("DS456", "Cybersecurity Techniques", "Elective", "4-0-0-0-4", 4, "Dr. Rajendra Hegadi"),
    @("EC456", "Reinforcement Learning", "Elective", "4-0-0-0-4", 4, "Dr. Divyajyothi"),
    @("DS401", "Health Care Data Analytics", "Elective", "4-0-0-0-4", 4, "Dr. Girish G N")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $summary.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

# Restore the originally active sheet/selection.
[void]$secA.Activate()
[void]$secA.Range("A1").Select()
